# Applies the "Added new achievements and translations" commit to Translations.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "Arkusz1" to "Translations"
$ws.Name = "Translations"

# 2. Insert two new rows before the current row 31 (Select / Confirm pair)
#    Existing rows 31-36 shift down to 33-38.
$ws.Rows.Item(31).Insert()
$ws.Rows.Item(31).Insert()

# Fill column-by-column to reproduce the original authoring order
# (English first, then the Key, then the Polish translation).
$ws.Range("C31").Value = "Select"
$ws.Range("C32").Value = "Confirm"

$ws.Range("A31").Value = "SELECT"
$ws.Range("A32").Value = "CONFIRM"

$ws.Range("B31").Value = "Wska" + [char]0x017C
$ws.Range("B32").Value = "Potwierd" + [char]0x017A

# 3. Append four new achievement rows at the bottom (rows 39-42)
$ws.Range("A39").Value = "SOLVEQUIZUNDER120"
$ws.Range("A40").Value = "SOLVEQUIZUNDER60"
$ws.Range("A41").Value = "SOLVEQUIZUNDER30"

$ws.Range("B39").Value = "Rozwi" + [char]0x0105 + [char]0x017C + " dowolny quiz w mniej ni" + [char]0x017C + " 2 minuty"
$ws.Range("B41").Value = "Rozwi" + [char]0x0105 + [char]0x017C + " dowolny quiz w mniej ni" + [char]0x017C + " 30 sekund"
$ws.Range("B40").Value = "Rozwi" + [char]0x0105 + [char]0x017C + " dowolny quiz w mniej ni" + [char]0x017C + " 1 minut" + [char]0x0119

$ws.Range("C39").Value = "Solve any quiz under 2 minutes"
$ws.Range("C40").Value = "Solve any quiz under 1 minute"
$ws.Range("C41").Value = "Solve any quiz under 30 seconds"

$ws.Range("A42").Value = "SOLVEQUIZUNDERLAST5"
$ws.Range("B42").Value = "Rozwi" + [char]0x0105 + [char]0x017C + " dowolny quiz w ostanich 5 sekundach"
$ws.Range("C42").Value = "Solve any quiz in last 5 seconds"

# 4. Restore the selection state shown in the final workbook
$ws.Range("C42").Select()
